$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Test Suite")

# Change the "Runmode" values in C2 and C3 from "Y" to "N"
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"

# Update the selected cell/range to match the saved view state
$ws.Range("B7").Select()
